$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recovery project list-Table")

# Update the Taxon text for the "Temminck's Ground Pangolin Reintroduction" row
# to wrap the scientific (Latin) name in asterisks, matching the markdown-style
# italics used for the other species entries in this column.
$ws.Range("D6").Value = "Ground pangolin (*Smutsia temminckii*)"
